# "abstracted methods in homepage"
# - Remove the "Spain, 1000-1200: Art at the Frontiers of Faith" exhibition row
#   from the Exhibitions sheet (rows below shift up).
# - Add two new sheets after "Exhibitions": "Explore" (a list of homepage link
#   URLs) and "Rules" (a list of visitor-guideline messages). "Rules" becomes
#   the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Exhibitions sheet: drop the "Spain, 1000-1200..." row entirely.
#    It is row 9 (A9) in the sheet as it stands before this edit.
# ---------------------------------------------------------------------------
$exhibitions = $wb.Worksheets.Item("Exhibitions")
[void]$exhibitions.Rows.Item(9).Delete()
[void]$exhibitions.Range("A2:A24").Select()

# ---------------------------------------------------------------------------
# 2. Add "Explore" sheet right after "Exhibitions" with homepage link URLs.
# ---------------------------------------------------------------------------
$explore = $wb.Worksheets.Add($null, $exhibitions)
$explore.Name = "Explore"

$exploreLinks = @(
    "https://www.metmuseum.org/perspectives",
    "https://www.metmuseum.org/150",
    "https://www.metmuseum.org/art/object-package?pkgids=725",
    "https://www.metmuseum.org/join-and-give/travel-with-the-met",
    "https://www.metmuseum.org/events/programs/virtual-events",
    "https://www.metmuseum.org/learn",
    "https://www.metmuseum.org/art/metpublications",
    "https://www.metmuseum.org/visit/group-visits",
    "https://www.metmuseum.org/visit/audio-content",
    "https://www.metmuseum.org/about-the-met/conservation-and-scientific-research",
    "https://www.metmuseum.org/about-the-met/collection-areas"
)

$row = 2
foreach ($link in $exploreLinks) {
    $explore.Cells.Item($row, 1).Value = $link
    $row++
}
[void]$explore.Range("A2:A12").Select()

# ---------------------------------------------------------------------------
# 3. Add "Rules" sheet right after "Explore" with visitor-guideline text.
# ---------------------------------------------------------------------------
$rules = $wb.Worksheets.Add($null, $explore)
$rules.Name = "Rules"

$rulesText = @(
    "Face coverings are required for all visitors age two and older, even if you are vaccinated.",
    "Maintain physical distancing. Keep at least six feet from others. ",
    "Report symptoms of illness.",
    "Wash hands and use hand sanitizer regularly during your visit.",
    "Select food and beverage options are available.",
    "Coat check is temporarily closed. No large bags. ",
    "Please follow directions from Met staff. ",
    "Please do not touch the art."
)

$row = 2
foreach ($line in $rulesText) {
    $rules.Cells.Item($row, 1).Value = $line
    $row++
}
[void]$rules.Range("A2:A9").Select()

# "Rules" ends up the active / selected sheet (tabSelected + workbookView activeTab).
[void]$rules.Activate()
